$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 801.2
$ws.Range("I12").Value = 744.8570999999999
$ws.Range("J12").Value = 932.6667
$ws.Range("K12").Value = 744.8570999999999
$ws.Range("L12").Value = 932.6667
$ws.Range("M12").Value = -574.8570999999999
$ws.Range("N12").Value = -1272.6667

# Row 64
$ws.Range("H64").Value = 3666.6667
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -4496

# Row 67
$ws.Range("H67").Value = 3666.6667
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -5716

# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents() | Out-Null

# Row 137
$ws.Range("H137").Value = 1700.3889
$ws.Range("I137").Value = 1286
$ws.Range("K137").Value = 3858
$ws.Range("M137").Value = -1308

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 584
$ws.Range("I2").Value = 563.4286
$ws.Range("K2").Value = 563.4286
$ws.Range("M2").Value = -450.4286

# Row 32
$ws.Range("H32").Value = 4016.173
$ws.Range("I32").Value = 2514.2563
$ws.Range("K32").Value = 2514.2563
$ws.Range("M32").Value = -2227.2563

# Row 45
$ws.Range("H45").Value = 1927.5
$ws.Range("I45").Value = 1043.3334
$ws.Range("K45").Value = 1043.3334
$ws.Range("M45").Value = -666.3334

# Row 97
$ws.Range("H97").Value = 549.5
$ws.Range("I97").Value = 549.5
$ws.Range("K97").Value = 549.5
$ws.Range("M97").Value = -53.5

# Row 116
$ws.Range("H116").Value = 584
$ws.Range("I116").Value = 563.4286
$ws.Range("K116").Value = 563.4286
$ws.Range("M116").Value = 1730.5714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 584
$ws.Range("I3").Value = 563.4286
$ws.Range("K3").Value = 563.4286
$ws.Range("M3").Value = -449.4286

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents() | Out-Null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2818.3928
$ws.Range("I31").Value = 1016.44446
$ws.Range("K31").Value = 1016.44446
$ws.Range("M31").Value = -721.44446

# Row 34
$ws.Range("H34").Value = 2818.3928
$ws.Range("I34").Value = 1016.44446
$ws.Range("K34").Value = 1016.44446
$ws.Range("M34").Value = -814.44446

# Row 58
$ws.Range("H58").Value = 1370.1578
$ws.Range("J58").Value = 1796.4286
$ws.Range("L58").Value = 1796.4286
$ws.Range("N58").Value = -2202.4286

# Row 70
$ws.Range("H70").Value = 29000
$ws.Range("J70").Value = 29000
$ws.Range("L70").Value = 29000
$ws.Range("N70").Value = -29630

# Row 73
$ws.Range("H73").Value = 29000
$ws.Range("J73").Value = 29000
$ws.Range("L73").Value = 29000
$ws.Range("N73").Value = -31184

# Row 99
$ws.Range("H99").Value = 1711.2858
$ws.Range("I99").Value = 1496.5
$ws.Range("K99").Value = 1496.5
$ws.Range("M99").Value = 1.5

# Row 122
$ws.Range("H122").Value = 1643.1724
$ws.Range("I122").Value = 1508.8235
$ws.Range("K122").Value = 4526.470499999999
$ws.Range("M122").Value = -2076.470499999999

# Row 126
$ws.Range("H126").Value = 1711.2858
$ws.Range("I126").Value = 1496.5
$ws.Range("K126").Value = 4489.5
$ws.Range("M126").Value = -2019.5

# Row 132
$ws.Range("H132").Value = 2316.3333
$ws.Range("I132").Value = 1240.3158
$ws.Range("J132").Value = 6405.2
$ws.Range("K132").Value = 3720.9474
$ws.Range("L132").Value = 19215.6
$ws.Range("M132").Value = -1190.9474
$ws.Range("N132").Value = -24275.6

# Row 136
$ws.Range("H136").Value = 1370.1578
$ws.Range("J136").Value = 1796.4286
$ws.Range("L136").Value = 5389.2858
$ws.Range("N136").Value = -10489.2858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 84
$ws.Range("J12").Value = 97.57143000000001
$ws.Range("L12").Value = 292.71429
$ws.Range("N12").Value = -638.71429

# Row 63
$ws.Range("H63").Value = 3506
$ws.Range("I63").Value = 3012
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 9036
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -8287
$ws.Range("N63").Value = -13498

# Row 66
$ws.Range("H66").Value = 3506
$ws.Range("I66").Value = 3012
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 27108
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -23364
$ws.Range("N66").Value = -43488

# Row 107
$ws.Range("H107").Value = 557.2778
$ws.Range("J107").Value = 557.2778
$ws.Range("L107").Value = 1671.8334
$ws.Range("N107").Value = -5511.8334

# Row 114
$ws.Range("H114").Value = 1972.8572
$ws.Range("I114").Value = 89.5
$ws.Range("K114").Value = 268.5
$ws.Range("M114").Value = 2985.5

# Row 131
$ws.Range("H131").Value = 10221107
$ws.Range("I131").Value = 71429150
$ws.Range("J131").Value = 19765.191
$ws.Range("K131").Value = 214287450
$ws.Range("L131").Value = 59295.573
$ws.Range("M131").Value = -214282410
$ws.Range("N131").Value = -69375.573

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5142.2856
$ws.Range("I70").Value = 5399.2
$ws.Range("K70").Value = 5399.2
$ws.Range("M70").Value = -5129.2

# Row 73
$ws.Range("H73").Value = 5142.2856
$ws.Range("I73").Value = 5399.2
$ws.Range("K73").Value = 5399.2
$ws.Range("M73").Value = -4463.2

# Row 86
$ws.Range("H86").Value = 30000
$ws.Range("J86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("N86").Value = -32372

# Row 89
$ws.Range("H89").Value = 30000
$ws.Range("J89").Value = 30000
$ws.Range("L89").Value = 90000
$ws.Range("N89").Value = -101856

# Row 122
$ws.Range("H122").Value = 1883
$ws.Range("I122").Value = 1772.125
$ws.Range("J122").Value = 2060.4
$ws.Range("K122").Value = 5316.375
$ws.Range("L122").Value = 6181.200000000001
$ws.Range("M122").Value = -2866.375
$ws.Range("N122").Value = -11081.2

# Row 126
$ws.Range("H126").Value = 62068.766
$ws.Range("J126").Value = 335330
$ws.Range("L126").Value = 1005990
$ws.Range("N126").Value = -1010930

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3496.6667
$ws.Range("I68").Value = 3196
$ws.Range("K68").Value = 3196
$ws.Range("M68").Value = -2447

# Row 71
$ws.Range("H71").Value = 3496.6667
$ws.Range("I71").Value = 3196
$ws.Range("K71").Value = 15980
$ws.Range("M71").Value = -12236

# Row 81
$ws.Range("H81").Value = 100000
$ws.Range("J81").Value = 100000
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996

# Row 82
$ws.Range("H82").Value = 3975.4443
$ws.Range("J82").Value = 4980
$ws.Range("L82").Value = 4980
$ws.Range("N82").Value = -5702

# Row 84
$ws.Range("H84").Value = 100000
$ws.Range("J84").Value = 100000
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984

# Row 85
$ws.Range("H85").Value = 3975.4443
$ws.Range("J85").Value = 4980
$ws.Range("L85").Value = 4980
$ws.Range("N85").Value = -7476

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 6000
$ws.Range("K56").Value = 6000
$ws.Range("M56").Value = -5286

# Row 92
$ws.Range("H92").Value = 27782.5
$ws.Range("J92").Value = 27782.5
$ws.Range("L92").Value = 27782.5
$ws.Range("N92").Value = -32774.5

# Row 132
$ws.Range("H132").Value = 4009.7896
$ws.Range("I132").Value = 1400.1111
$ws.Range("J132").Value = 6358.5
$ws.Range("K132").Value = 4200.3333
$ws.Range("L132").Value = 19075.5
$ws.Range("M132").Value = -1670.3333
$ws.Range("N132").Value = -24135.5

# Row 136
$ws.Range("H136").Value = 3090.7896
$ws.Range("I136").Value = 2273.4
$ws.Range("K136").Value = 6820.200000000001
$ws.Range("M136").Value = -4270.200000000001
